$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    334 = 0.9599317154288292
    335 = 0.9455632781237364
    336 = 0.7641119734942913
    337 = 1.164684123396873
    338 = 1.156321626976132
    339 = 0.9391011016070843
    340 = 0.6331674814224243
    341 = 0.5372912478446961
    342 = 0.7946128785610199
    343 = 0.7913253098726273
    344 = 0.5530621093511582
    345 = 0.4222658997774124
    346 = 0.2527838963270188
    347 = 0.4344980749487877
    348 = 0.1481849277019501
    349 = 0.4633985915780067
    350 = 0.5122897350788117
    351 = 0.2426398795843124
    352 = 1.27256611071527
    353 = 1.022898229211569
    354 = 0.1922542336583138
    355 = 0.6728441748023033
    356 = 0.6459268774092197
    357 = 0.5305473661422729
    358 = 0.5466447226703167
    359 = 0.3594466164708138
    360 = 0.2048006477952004
    361 = 0.1636517548561096
    362 = 0.8283556419610977
    363 = 0.549017400443554
    364 = 0.0838139119744301
    365 = 0.5645945847034455
    366 = 0.4806778948009014
    367 = 0.2398127809166908
    368 = 0.4772502493858338
    369 = 0.3967638349533081
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
